# Update the player table (Oyuncu Adı / Pozisyon / Takım) on the active sheet.
# Rows 2-18 get new Name / Position / Team values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox",        "PG",       "Sacramento Kings"),
    @("Kelly Oubre Jr.",     "SG,SF",    "Philadelphia 76ers"),
    @("Bennedict Mathurin",  "SG,SF",    "Indiana Pacers"),
    @("DeMar DeRozan",       "SF,PF",    "Sacramento Kings"),
    @("Guerschon Yabusele",  "PF,C",     "Philadelphia 76ers"),
    @("Nikola Vucevic",      "PF,C",     "Chicago Bulls"),
    @("Brook Lopez",         "C",        "Milwaukee Bucks"),
    @("Tyler Herro",         "PG,SG",    "Miami Heat"),
    @("Ja Morant",           "PG",       "Memphis Grizzlies"),
    @("Amen Thompson",       "SG,SF",    "Houston Rockets"),
    @("Mikal Bridges",       "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes",      "SG,SF,PF", "Toronto Raptors"),
    @("Luka Doncic",         "PG,SG",    "Dallas Mavericks"),
    @("Evan Mobley",         "PF,C",     "Cleveland Cavaliers"),
    @("Santi Aldama",        "PF,C",     "Memphis Grizzlies"),
    @("Miles Bridges",       "SF,PF",    "Charlotte Hornets"),
    @("Josh Giddey",         "PG,SG,SF", "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
